# Applies the "Atualizado por script em 05-11-2023 08:45" update to the
# san-marino / campionato-sammarinese 2023-2024 sheet:
#  1) Three groups of rows had their match-data columns (F:V) rotated
#     among themselves (the site re-sorted/re-fetched the fixtures),
#     while columns A:E (index/pais/torneio/temporada/data_partida)
#     stayed put.
#  2) Five brand-new match rows (51-55) were appended at the bottom.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Rotate F:V among rows 26, 27, 29  (new26=old29, new27=old26, new29=old27)
# ---------------------------------------------------------------------
$vals26 = $ws.Range("F26:V26").Value2
$vals27 = $ws.Range("F27:V27").Value2
$vals29 = $ws.Range("F29:V29").Value2

$ws.Range("F26:V26").Value2 = $vals29
$ws.Range("F27:V27").Value2 = $vals26
$ws.Range("F29:V29").Value2 = $vals27

# ---------------------------------------------------------------------
# Rotate F:V among rows 44, 45, 46  (new44=old46, new45=old44, new46=old45)
# ---------------------------------------------------------------------
$vals44 = $ws.Range("F44:V44").Value2
$vals45 = $ws.Range("F45:V45").Value2
$vals46 = $ws.Range("F46:V46").Value2

$ws.Range("F44:V44").Value2 = $vals46
$ws.Range("F45:V45").Value2 = $vals44
$ws.Range("F46:V46").Value2 = $vals45

# ---------------------------------------------------------------------
# Rotate F:V among rows 48, 49, 50  (new48=old49, new49=old50, new50=old48)
# ---------------------------------------------------------------------
$vals48 = $ws.Range("F48:V48").Value2
$vals49 = $ws.Range("F49:V49").Value2
$vals50 = $ws.Range("F50:V50").Value2

$ws.Range("F48:V48").Value2 = $vals49
$ws.Range("F49:V49").Value2 = $vals50
$ws.Range("F50:V50").Value2 = $vals48

# ---------------------------------------------------------------------
# 2) Append five new rows (51-55) with the same look & feel (formatting)
#    as the existing data rows, then fill in their values.
# ---------------------------------------------------------------------
$ws.Range("A50:V50").Copy()
$ws.Range("A51:V55").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Cells.Item(51, 1).Value2 = 50
$ws.Cells.Item(51, 2).Value2 = "san-marino"
$ws.Cells.Item(51, 3).Value2 = "campionato-sammarinese"
$ws.Cells.Item(51, 4).Value2 = "2023-2024"
$ws.Cells.Item(51, 5).Value2 = 45234.625
$ws.Cells.Item(51, 6).Value2 = "Faetano"
$ws.Cells.Item(51, 7).Value2 = 4
$ws.Cells.Item(51, 8).Value2 = "Folgore"
$ws.Cells.Item(51, 9).Value2 = 1
$ws.Cells.Item(51, 10).Value2 = 3.49
$ws.Cells.Item(51, 11).Value2 = "03/11/2023 03:12"
$ws.Cells.Item(51, 12).Value2 = 3.69
$ws.Cells.Item(51, 13).Value2 = "04/11/2023 14:05"
$ws.Cells.Item(51, 14).Value2 = 3.48
$ws.Cells.Item(51, 15).Value2 = "03/11/2023 03:12"
$ws.Cells.Item(51, 16).Value2 = 3.85
$ws.Cells.Item(51, 17).Value2 = "04/11/2023 14:16"
$ws.Cells.Item(51, 18).Value2 = 1.74
$ws.Cells.Item(51, 19).Value2 = "03/11/2023 03:12"
$ws.Cells.Item(51, 20).Value2 = 1.74
$ws.Cells.Item(51, 21).Value2 = "04/11/2023 14:08"
$ws.Cells.Item(51, 22).Value2 = "https://www.betexplorer.com/football/san-marino/campionato-sammarinese/sc-faetano-folgore/thsTZ6XD/"

$ws.Cells.Item(52, 1).Value2 = 51
$ws.Cells.Item(52, 2).Value2 = "san-marino"
$ws.Cells.Item(52, 3).Value2 = "campionato-sammarinese"
$ws.Cells.Item(52, 4).Value2 = "2023-2024"
$ws.Cells.Item(52, 5).Value2 = 45234.625
$ws.Cells.Item(52, 6).Value2 = "Cosmos"
$ws.Cells.Item(52, 7).Value2 = 1
$ws.Cells.Item(52, 8).Value2 = "Tre Penne"
$ws.Cells.Item(52, 9).Value2 = 1
$ws.Cells.Item(52, 10).Value2 = 2.6
$ws.Cells.Item(52, 11).Value2 = "03/11/2023 03:12"
$ws.Cells.Item(52, 12).Value2 = 3.21
$ws.Cells.Item(52, 13).Value2 = "04/11/2023 14:16"
$ws.Cells.Item(52, 14).Value2 = 3.01
$ws.Cells.Item(52, 15).Value2 = "03/11/2023 03:12"
$ws.Cells.Item(52, 16).Value2 = 3.15
$ws.Cells.Item(52, 17).Value2 = "04/11/2023 14:16"
$ws.Cells.Item(52, 18).Value2 = 2.32
$ws.Cells.Item(52, 19).Value2 = "03/11/2023 03:12"
$ws.Cells.Item(52, 20).Value2 = 2.1
$ws.Cells.Item(52, 21).Value2 = "04/11/2023 14:06"
$ws.Cells.Item(52, 22).Value2 = "https://www.betexplorer.com/football/san-marino/campionato-sammarinese/sp-cosmos-tre-penne/hU9dIv2t/"

$ws.Cells.Item(53, 1).Value2 = 52
$ws.Cells.Item(53, 2).Value2 = "san-marino"
$ws.Cells.Item(53, 3).Value2 = "campionato-sammarinese"
$ws.Cells.Item(53, 4).Value2 = "2023-2024"
$ws.Cells.Item(53, 5).Value2 = 45234.625
$ws.Cells.Item(53, 6).Value2 = "Libertas"
$ws.Cells.Item(53, 7).Value2 = 2
$ws.Cells.Item(53, 8).Value2 = "Cailungo"
$ws.Cells.Item(53, 9).Value2 = 1
$ws.Cells.Item(53, 10).Value2 = 1.43
$ws.Cells.Item(53, 11).Value2 = "03/11/2023 03:12"
$ws.Cells.Item(53, 12).Value2 = 1.51
$ws.Cells.Item(53, 13).Value2 = "04/11/2023 14:05"
$ws.Cells.Item(53, 14).Value2 = 4.12
$ws.Cells.Item(53, 15).Value2 = "03/11/2023 03:12"
$ws.Cells.Item(53, 16).Value2 = 4.42
$ws.Cells.Item(53, 17).Value2 = "04/11/2023 14:05"
$ws.Cells.Item(53, 18).Value2 = 4.83
$ws.Cells.Item(53, 19).Value2 = "03/11/2023 03:12"
$ws.Cells.Item(53, 20).Value2 = 4.61
$ws.Cells.Item(53, 21).Value2 = "04/11/2023 14:05"
$ws.Cells.Item(53, 22).Value2 = "https://www.betexplorer.com/football/san-marino/campionato-sammarinese/ac-libertas-cailungo/rJ29Fxna/"

$ws.Cells.Item(54, 1).Value2 = 53
$ws.Cells.Item(54, 2).Value2 = "san-marino"
$ws.Cells.Item(54, 3).Value2 = "campionato-sammarinese"
$ws.Cells.Item(54, 4).Value2 = "2023-2024"
$ws.Cells.Item(54, 5).Value2 = 45234.625
$ws.Cells.Item(54, 6).Value2 = "San Giovanni"
$ws.Cells.Item(54, 7).Value2 = 0
$ws.Cells.Item(54, 8).Value2 = "La Fiorita"
$ws.Cells.Item(54, 9).Value2 = 3
$ws.Cells.Item(54, 10).Value2 = 8.47
$ws.Cells.Item(54, 11).Value2 = "03/11/2023 03:12"
$ws.Cells.Item(54, 12).Value2 = 8.97
$ws.Cells.Item(54, 13).Value2 = "04/11/2023 14:18"
$ws.Cells.Item(54, 14).Value2 = 4.93
$ws.Cells.Item(54, 15).Value2 = "03/11/2023 03:12"
$ws.Cells.Item(54, 16).Value2 = 4.58
$ws.Cells.Item(54, 17).Value2 = "04/11/2023 14:18"
$ws.Cells.Item(54, 18).Value2 = 1.21
$ws.Cells.Item(54, 19).Value2 = "03/11/2023 03:12"
$ws.Cells.Item(54, 20).Value2 = 1.29
$ws.Cells.Item(54, 21).Value2 = "04/11/2023 14:18"
$ws.Cells.Item(54, 22).Value2 = "https://www.betexplorer.com/football/san-marino/campionato-sammarinese/san-giovanni-la-fiorita/SlgoBNfD/"

$ws.Cells.Item(55, 1).Value2 = 54
$ws.Cells.Item(55, 2).Value2 = "san-marino"
$ws.Cells.Item(55, 3).Value2 = "campionato-sammarinese"
$ws.Cells.Item(55, 4).Value2 = "2023-2024"
$ws.Cells.Item(55, 5).Value2 = 45234.76041666666
$ws.Cells.Item(55, 6).Value2 = "Domagnano"
$ws.Cells.Item(55, 7).Value2 = 0
$ws.Cells.Item(55, 8).Value2 = "Virtus"
$ws.Cells.Item(55, 9).Value2 = 1
$ws.Cells.Item(55, 10).Value2 = 6.07
$ws.Cells.Item(55, 11).Value2 = "03/11/2023 06:42"
$ws.Cells.Item(55, 12).Value2 = 11.47
$ws.Cells.Item(55, 13).Value2 = "04/11/2023 17:49"
$ws.Cells.Item(55, 14).Value2 = 4.33
$ws.Cells.Item(55, 15).Value2 = "03/11/2023 06:42"
$ws.Cells.Item(55, 16).Value2 = 5.77
$ws.Cells.Item(55, 17).Value2 = "04/11/2023 17:49"
$ws.Cells.Item(55, 18).Value2 = 1.32
$ws.Cells.Item(55, 19).Value2 = "03/11/2023 06:42"
$ws.Cells.Item(55, 20).Value2 = 1.18
$ws.Cells.Item(55, 21).Value2 = "04/11/2023 17:49"
$ws.Cells.Item(55, 22).Value2 = "https://www.betexplorer.com/football/san-marino/campionato-sammarinese/sp-domagnano-virtus/zga1HbHn/"

Write-Host "Edit complete: rotated rows 26/27/29, 44/45/46, 48/49/50 and appended rows 51-55"
